# Update "想去人数" (want-to-go count) figures for a handful of events
# across the 展览 (sheet1) and 全部类型 (sheet4) sheets.

$wb = $excel.ActiveWorkbook

# 展览 sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 5384
$wsExpo.Range("F6").Value = 70
$wsExpo.Range("F7").Value = 1
$wsExpo.Range("F9").Value = 517

# 全部类型 sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5384
$wsAll.Range("F7").Value = 70
$wsAll.Range("F8").Value = 1
$wsAll.Range("F11").Value = 517
